$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(296).Insert()

$ws.Range("A296").Value = 3
$ws.Range("B296").Value = "Femacal de La Calera"
$ws.Range("C296").Value = "Coquimbo"
$ws.Range("D296").Value = 44754
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = 100112009
$ws.Range("G296").Value = "Acelga"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 260
$ws.Range("K296").Value = 3000
$ws.Range("L296").Value = 3300
$ws.Range("M296").Value = 3162
$ws.Range("N296").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O296").Value = "Provincia de Quillota"
$ws.Range("P296").Value = 527
$ws.Range("Q296").Value = 6
$ws.Range("R296").Value = "Hortaliza"
